$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the last-changed date for every data row
# (rows 2 through 267). All of them move from 2023-09-17 (serial 45186)
# to 2023-09-19 (serial 45188).
$range = $ws.Range("C2:C267")
$range.Value = 45188
